$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: wrap a <w:body> fragment into a minimal single-part WordOpenXML
# package suitable for Range.InsertXML.
# ---------------------------------------------------------------------------
function New-WordPackageXml([string]$bodyFragment) {
    return '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyFragment + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ===========================================================================
# Edit 1: title paragraph - split "Breadth First Search)" into separate runs
# wrapped with spell-check proofErr markers (spellStart/spellEnd), keeping
# the same run formatting (color + size) throughout.
# ===========================================================================

$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleEnd = $titleRange.End

$rPrFrag = '<w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/><w:sz w:val="44"/></w:rPr>'

$newTitleBody = '<w:p><w:pPr>' + $rPrFrag + '</w:pPr>' +
    '<w:r>' + $rPrFrag + '<w:t>Algo</w:t></w:r>' +
    '<w:r>' + $rPrFrag + '<w:t xml:space="preserve">ritmo de B&#250;squeda en Anchura (BFS: </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPrFrag + '<w:t>Breadth</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPrFrag + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPrFrag + '<w:t>First</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPrFrag + '<w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r>' + $rPrFrag + '<w:t>Search</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r>' + $rPrFrag + '<w:t>)</w:t></w:r>' +
    '</w:p>'

# Inserting a <w:p>-wrapped fragment at the very end of a paragraph (just
# before its paragraph mark) creates a brand-new paragraph right after it.
$insertionPoint = $d.Range($titleEnd - 1, $titleEnd - 1)
$insertionPoint.InsertXML((New-WordPackageXml $newTitleBody))

# The original (old-formatted) title paragraph is now redundant -- remove it,
# leaving our freshly inserted paragraph in its place as paragraph 1.
$oldTitlePara = $d.Paragraphs(1)
$oldTitlePara.Range.Delete()

# ===========================================================================
# Edit 2: append a new run after "Cantidad de Ramas del Recorrido"
# ===========================================================================

$target = $d.Content
$found = $target.Find.Execute("Cantidad de Ramas del Recorrido", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0)

$targetPara = $target.Paragraphs(1)
$paraRange = $targetPara.Range
$paraStart = $paraRange.Start

$newBody = '<w:p><w:r><w:t>Cantidad de Ramas del Recorrido</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> &#8211; Nodo Inicial (Ra&#237;z del &#193;rbol del Recorrido)</w:t></w:r></w:p>'

$startPoint = $d.Range($paraStart, $paraStart)
$startPoint.InsertXML((New-WordPackageXml $newBody))

# Inserting a <w:p>-wrapped fragment at the very start of an existing
# paragraph merges into it (no new paragraph break), but leaves the
# paragraph's original run text duplicated at the tail -- remove that.
$updatedPara = $d.Paragraphs(8)
$updatedRange = $updatedPara.Range
$afterFirst = $d.Range($updatedRange.Start, $updatedRange.Start)
$afterFirst.Find.Execute("Cantidad de Ramas del Recorrido", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterFirst.Collapse(0)

$tail = $d.Range($afterFirst.End, $updatedRange.End)
$tail.Find.Execute("Cantidad de Ramas del Recorrido", $true, $false,
    $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$tail.Delete()

Write-Output "done"
